$d = $word.ActiveDocument

# Sanity-check the starting shape of the document before mutating it.
$startCount = $d.Paragraphs.Count
if ($startCount -ne 20) {
    Write-Output ("WARNING: expected 20 paragraphs before edit, found " + $startCount)
}

# Paragraph 1: date line gets updated + a manual line break + new title line added
$d.Paragraphs.Item(1).Range.Text = "המאמר היומי של מייק: 07.07.25`vProcedural Knowledge in Pretraining Drives Reasoning in Large Language Models"

$d.Paragraphs.Item(2).Range.Text = "מודלי שפה גדולים ממשיכים להדהים אותנו ביכולותיהם האדירות, אך שאלה מטרידה נותרה בעינה: האם הם באמת `"מבינים`", או שהם פשוט תוכים מתוחכמים המשננים את דאטה האימון שלהם? המאמר המסוקר מציע פרספקטיבה חדשה, החורגת ממגבלות ההפרדה המסורתיות של דאטהסטי אימון וטסט כדי לחקור כיצד LLMs לומדים ״להסיק מסקנות״ מדאטה של האימון המקדים שלהם. (pretraining)."
$d.Paragraphs.Item(3).Range.Text = "הגודל העצום של דאטהסט לאימון מקדים של LLMs הקשה היסטורית על הבחנה האם ביצועי מודל במשימה נובעים מהכללה אמיתית או משינון בלבד של דוגמאות שנתקלו בהן בעבר. נהמחברים מתמודדים עם זה על ידי שימוש בפונקציות השפעה, טכניקה מעולם הסטטיסטיקה, כדי לזהות אילו מסמכי אימון מקדים ספציפיים משפיעים על פלט המודל עבור שאילתות נתונות. גישתם חדשנית בהתמקדותה בהשפעת דאטה האימון המקדים במקום בפרשנות בלבד של משקלי ואקטיבציות המודל, ומספקת זווית ייחודית לתהליך הלמידה."
$d.Paragraphs.Item(4).Range.Text = "הגילוי המשמעותי ביותר הוא שעבור משימות הנמקה (במיוחד, בעיות מתמטיות כמו אריתמטיקה, חישוב שיפועים ופתרון משוואות לינאריות), השפעת מסמכי האימון המקדים מתואמת מאוד בין שאילתות שונות באותה משימה. משמעות הדבר היא שמסמך המשפיע על חישוב שיפוע אחד צפוי להשפיע גם על חישוב אחר, גם עם מספרים שונים. זה מצביע בחזקה על כך ש-LLMs אינם רק מאחזרים תשובות ספציפיות אלא מחלצים ומיישמים ידע פרוצדורלי שלבי `"איך לעשות`" או אלגוריתמים מהדאטה. זה עומד בניגוד חד לשאילתות עובדתיות, שבהן ההשפעה ספציפית מאוד לכל שאלה, מה שמצביע על אחזור ישיר יותר של עובדות משוננות."
$d.Paragraphs.Item(5).Range.Text = "המחקר מצא כי גודל ההשפעה ממסמכים בודדים נמוך באופן כללי עבור שאלות חשיבה בהשוואה לשאלות עובדתיות. יתר על כן, קבוצת המסמכים המשפיעים על חשיבה פחות `"ספציפית`" ויותר כללית. משמעות הדבר היא שעבור חשיבה, LLMs שואבים ממערך ידע רחב ומפוזר יותר, ומסתמכים פחות על מסמך בודד כלשהו. ממצא זה תומך ברעיון של אסטרטגיית למידה מוכללת יותר עבור חשיבה, שבה המודל מסנתז מידע ממקורות רבים במקום לאתר כמה מקורות רלוונטיים במיוחד. ההשפעה בולטת אף יותר במודלים גדולים יותר, מה שמצביע על יעילות נתונים גבוהה יותר בהכללה."
$d.Paragraphs.Item(6).Range.Text = "באופן מסקרן, בעוד שתשובות לשאלות עובדתיות מופיעות לעתים קרובות ב-0.01% העליונים של מסמכי האימון המקדים המשפיעים, זה כמעט אף פעם לא המקרה עבור שאלות חשיבה. גם כאשר שלבי חשיבה ביניים או תשובות מלאות קיימים במערך הנתונים הרחב יותר של האימון המקדים, הם מופיעים לעתים רחוקות כבעלי השפעה רבה על שאילתות חשיבה. זה מחזק עוד יותר את הרעיון ש-LLMs אינם פשוט `"מאחזרים`" את הפתרון לבעיית חשיבה אלא מיישמים פרוצדורות נלמדות."
$d.Paragraphs.Item(7).Range.Text = "המחקר מדגיש את התפקיד המשמעותי של קוד בהנעת יכולות חשיבה.דאטהסטים הקשורים לקוד (כמו StackExchange) נמצאו כבעלי ייצוג יתר משמעותי בקרב המסמכים המשפיעים ביותר על שאילתות חשיבה, הרבה מעבר לשיעורם בהתפלגות האימון הכוללת. זה מצביע על כך שקוד, עם המבנה הלוגי והפרוצדוראלי הטבוע בו, משמש כמקור עשיר עבור LLMs ללמוד אסטרטגיות חשיבה ניתנות להכללה. ממצא זה פותח אפיקים חדשים לאופטימיזציה של הרכב נתוני האימון המקדים כדי לשפר את החשיבה."
$d.Paragraphs.Item(8).Range.Text = "ממצאים אלה מאתגרים את התפיסה הפשטנית של LLMs כ`"תוכי סטוכסטי`", לפחות בכל הנוגע ליכולות החשיבה שלהם. במקום פשוט לחזור על מידע, נראה שהמודלים לומדים פרוצדורות מופשטות ומיישמים אותן לבעיות חדשות. הכללה פרוצדורלית זו היא צעד קריטי לקראת בינה מלאכותית חזקה ואמיתית יותר."
$d.Paragraphs.Item(9).Range.Text = "ההשלכות לפיתוח LLM עתידיות :"
$d.Paragraphs.Item(10).Range.Text = "במקום לנסות לכסות כל מקרה אפשרי של בעיה, אסטרטגיות אימון מקדים יכולות להתמקד בדאטה באיכות גבוהה המדגימים במפורש פרוצדורות ומתודולוגיות לפתרון בעיות על פני משימות חשיבה מגוונות."
$d.Paragraphs.Item(11).Range.Text = "ההשפעה העצומה של קוד מצביעה על כך שהגדלת נוכחותם או אוצרותם באופן ספציפי עבור התוכן הפרוצדורלי שלהם יכולה להיות דרך יעילה ביותר להגביר את חשיבת ה-LLM."
$d.Paragraphs.Item(12).Range.Text = "ההבנה שהנמקה אינה רק אחזור מאפשרת לנו לתכנן אמות מידה ומדדי הערכה טובים יותר הבודקים באמת את יכולת המודל להכליל וליישם פרוצדורות נלמדות."

# Remove the now-obsolete paragraphs 13-19 (old scheming-paper body + Heading3 section)
$pStart = $d.Paragraphs.Item(13)
$pEnd = $d.Paragraphs.Item(19)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rng.Delete()

# Update the arxiv link at the end of the document
$d.Content.Find.Execute("https://arxiv.org/abs/2412.04984", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/abs/2411.12580", 2) | Out-Null

$endCount = $d.Paragraphs.Count
if ($endCount -ne 13) {
    Write-Output ("WARNING: expected 13 paragraphs after edit, found " + $endCount)
}
Write-Output ("Final paragraph count=" + $endCount)
